$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "33/52"
$ws.Range("H3").Value = "52/52"
$ws.Range("H4").Value = "42/52"
$ws.Range("H5").Value = "35/52"
$ws.Range("H6").Value = "45/52"
$ws.Range("G7").Value = "admin@admin.com, System"
$ws.Range("H7").Value = "35/52"
$ws.Range("H8").Value = "44/52"
$ws.Range("H9").Value = "21/52"
$ws.Range("H10").Value = "27/52"
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "71.2%"
$ws.Range("L10").NumberFormat = "general"
$ws.Range("H11").Value = "34/52"
$ws.Range("H12").Value = "31/52"
$ws.Range("H13").Value = "36/52"
$ws.Range("H14").Value = "36/52"
$ws.Range("H15").Value = "43/52"
$ws.Range("M15").Value = 52
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "70.0%"
$ws.Range("S15").NumberFormat = "general"
$ws.Range("H16").Value = "42/52"
$ws.Range("H17").Value = "36/52"
$ws.Range("H18").Value = "38/52"
$ws.Range("H19").Value = "35/52"
$ws.Range("M19").Value = 56
$ws.Range("S19").NumberFormat = "@"
$ws.Range("S19").Value = "73.0%"
$ws.Range("S19").NumberFormat = "general"
$ws.Range("H20").Value = "42/52"
$ws.Range("H21").Value = "42/52"
$ws.Range("H22").Value = "31/52"
$ws.Range("H23").Value = "41/52"
$ws.Range("H24").Value = "27/52"
$ws.Range("H25").Value = "26/52"
$ws.Range("H26").Value = "0/52"
$ws.Range("H27").Value = "0/52"
$ws.Range("H28").Value = "0/52"
$ws.Range("G34").Value = "admin@admin.com, System"
$ws.Range("G61").Value = "admin@admin.com, System"
$ws.Range("H109").Value = "30/56"
$ws.Range("H110").Value = "38/56"
$ws.Range("H111").Value = "55/56"
$ws.Range("H112").Value = "54/56"
$ws.Range("H113").Value = "26/56"
$ws.Range("H114").Value = "37/56"
$ws.Range("H115").Value = "39/56"
$ws.Range("H116").Value = "54/56"
$ws.Range("H117").Value = "54/56"
$ws.Range("H118").Value = "47/56"
$ws.Range("H119").Value = "42/56"
$ws.Range("H120").Value = "45/56"
$ws.Range("H121").Value = "43/56"
$ws.Range("H122").Value = "41/56"
$ws.Range("H123").Value = "36/56"
$ws.Range("H124").Value = "38/56"
$ws.Range("H125").Value = "45/56"
$ws.Range("H126").Value = "36/56"
$ws.Range("H127").Value = "33/56"
$ws.Range("H128").Value = "42/56"
$ws.Range("H129").Value = "29/56"
$ws.Range("H130").Value = "35/56"
$ws.Range("H131").Value = "0/56"
$ws.Range("H132").Value = "0/56"
$ws.Range("H133").Value = "0/56"
$ws.Range("H134").Value = "0/56"
